$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend table with a new "2020" column (K), matching formatting of the
# existing "2019" column (J) cell by cell.

# Row 3: empty bottom-border cell
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# Row 4: header year 2020
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2020

# Row 6: Mammals value
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 5.9

# Row 7: Birds value
$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("K7").Value = 1.5

# Row 8: Amphibians and Reptiles - no data ("-")
$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = "-"

# Restore the active selection recorded in the saved workbook
[void]$ws.Range("L16").Select()
